$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Customer detail section
$ws.Range("C2").Value = "SINTIA WOWOR"
$ws.Range("C3").Value = "`n                    16/06/2021  4:12PM"
$ws.Range("C4").Value = "'06111021006591      "
$ws.Range("C5").Value = "AIRMADIDI , 06/09/1993"
$ws.Range("C6").Value = "SINTIAWOWOR@GMAIL.COM"
$ws.Range("C7").Value = "'3,033,424.75"
$ws.Range("C8").Value = "'2,291,400.00"

# Asset view section
$ws.Range("C10").Value = "FIKLY FIRGIN MANTIRI"
$ws.Range("C11").Value = "SARONGSONG I LK V KEC AIRMADIDI KAB MINAHASA UTARA"
$ws.Range("C12").Value = "MH1KC8218HK101981"
$ws.Range("C13").Value = "KC82E1099191"
$ws.Range("C14").Value = "'2017"
$ws.Range("C15").Value = "'14916"
$ws.Range("C16").Value = "HITAM MERAH"
$ws.Range("C17").Value = "DB 3114 FH"
$ws.Range("C18").Value = "N01773265S"
$ws.Range("C19").Value = "FHDC001030Q"
